$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tbl = $shape.Table

$tbl.Rows.Item(2).Cells.Item(1).Shape.TextFrame.TextRange.Text = "A"
$tbl.Rows.Item(3).Cells.Item(1).Shape.TextFrame.TextRange.Text = "B"
$tbl.Rows.Item(4).Cells.Item(1).Shape.TextFrame.TextRange.Text = "C"
$tbl.Rows.Item(5).Cells.Item(1).Shape.TextFrame.TextRange.Text = "D"
